$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 1, 2, and 5 contents (these rows are removed entirely from the target)
$ws.Range("A1:C2").ClearContents()
$ws.Range("A5:C5").ClearContents()

# Row 4: remove "Zanahoria"/"2.0"/"KG" data and replace with "TOTAL"/"$ 0.0"
$ws.Range("A4").Value = "TOTAL"

# Leading apostrophe forces text interpretation so "$ 0.0" stays a literal
# string instead of being parsed as a currency number.
$ws.Range("B4").Value = "'$ 0.0"

$ws.Range("C4").ClearContents()
